$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '51.054.07'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '2.957.12'
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '381.61'
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("D6").Value = '102.16'
$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("D7").Value = '0.546'
$ws.Range("E7").Value = '  +1.85%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("D10").Value = '36.49'
$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("E11").Value = '  -0.65%  '

$ws.Range("D12").Value = '0.0850'
$ws.Range("E12").Value = '  +1.81%  '

$ws.Range("B13").Value = 'Uniswap'
$ws.Range("C13").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D13").Value = '12.46'
$ws.Range("E13").Value = '  +74.91%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.419.41'
$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '18.40'
$ws.Range("E15").Value = '  +2.49%  '

$ws.Range("D16").Value = '7.75'
$ws.Range("E16").Value = '  +5.43%  '

$ws.Range("D17").Value = '2.967.22'
$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("E18").Value = '  +3.56%  '

$ws.Range("D19").Value = '51.126.43'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = '3.08'
$ws.Range("E20").Value = '  -2.40%  '

$ws.Range("D21").Value = '12.40'
$ws.Range("E21").Value = '  -1.18%  '

$ws.Range("D22").Value = '0.0₃0958'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").Value = '3.35'
$ws.Range("E23").Value = '  +16.56%  '

$ws.Range("D24").Value = '269.85'
$ws.Range("E24").Value = '  +2.57%  '

$ws.Range("D25").Value = '69.75'
$ws.Range("E25").Value = '  +2.28%  '

$ws.Range("D26").Value = '7.93'
$ws.Range("E26").Value = '  -2.74%  '

$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").Value = '25.83'
$ws.Range("E29").Value = '  +0.57%  '

$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -11.52%  '

$ws.Range("E31").Value = '  -3.30%  '

$ws.Range("D32").Value = '10.48'
$ws.Range("E32").Value = '  +6.36%  '

$ws.Range("D33").Value = '51.19'
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("D34").Value = '34.30'

$ws.Range("E35").Value = '  +2.10%  '

$ws.Range("E36").Value = '  -4.24%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '3.25'
$ws.Range("E38").Value = '  +9.19%  '

$ws.Range("E39").Value = '  +2.15%  '

$ws.Range("D40").Value = '16.72'
$ws.Range("E40").Value = '  +1.55%  '

$ws.Range("E41").Value = '  +2.60%  '

$ws.Range("D42").Value = '2.51'
$ws.Range("E42").Value = '  -2.70%  '

$ws.Range("D43").Value = '124.46'
$ws.Range("E43").Value = '  +2.32%  '

$ws.Range("D44").Value = '21.77'
$ws.Range("E44").Value = '  +3.13%  '

$ws.Range("D45").Value = '3.55'
$ws.Range("E45").Value = '  +10.13%  '

$ws.Range("D46").Value = '2.084.33'
$ws.Range("E46").Value = '  +4.11%  '

$ws.Range("D47").Value = '2.03'
$ws.Range("E47").Value = '  -0.89%  '

$ws.Range("E48").Value = '  -0.31%  '

$ws.Range("E49").Value = '  -3.31%  '

$ws.Range("E50").Value = '  -6.70%  '

$ws.Range("E51").Value = '  +6.62%  '

# Restore default (General) formatting / no explicit style on column D
$ws.Range("D2:D51").ClearFormats()
